# First launch version 15 april
# Fix typo in the English COVID-19 scheduling notice: "EGG" -> "EEG"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B31").Value = "Due to the current regulations implemented in regard to the COVID-19 crisis, we will only schedule EEG experiments once it is permitted to do so"

# Update the view/selection state to match the saved workbook view
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("D39").Select()
